$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 4, pushing TIME_REQ..FINAL_STATE_REQ down.
# This splits the old "I_IN_REQ" entry (row 3) into two entries:
#   row 3 -> I_IN_1_REQ (exponent 1)
#   row 4 -> I_IN_2_REQ (exponent 2, newly inserted row)
$ws.Rows.Item(4).Insert()

# Rename row 3 from I_IN_REQ to I_IN_1_REQ
$ws.Range("A3").Value = "I_IN_1_REQ"

# Fill in the new row 4 with I_IN_2_REQ, mirroring the structure of row 3
$ws.Range("A4").Value = "I_IN_2_REQ"
$ws.Range("B4").Value = 2

# Renumber the exponent column (B) for all rows from the old row 4 (now row 5) downward
for ($r = 5; $r -le 14; $r++) {
    $ws.Range("B$r").Value = $r - 2
}

# Fix the MAX() range in G2 to include the new last row (14)
$ws.Range("G2").Formula = "=MAX(C2:C14)"

# Update/refresh formulas for C/D/E columns across rows 3:14 (row 4 was already set above)
for ($r = 3; $r -le 14; $r++) {
    $ws.Range("C$r").Formula = "=LEN(A$r)"
    $ws.Range("D$r").Formula = '=CONCATENATE(A' + $r + ', REPT(" ",$G$2+1-C' + $r + '))'
    $ws.Range("E$r").Formula = '=CONCATENATE(D' + $r + ',' + '" = (1 << "' + ',B' + $r + ',"), ")'
}

# Update selection to match the new active cell shown in the diff
$ws.Range("F14").Select()

$wb.Application.Calculate()
